$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'hamstring knee compression',
    'hamstring recovery',
    'hamstring stick',
    'hamstring stretch',
    'hamstring stretches',
    'hamstring support',
    'hamstring support for men',
    'hamstring support for women',
    'hamstring workout equipment',
    'hamstrings compression',
    'hamstrings stretch',
    'hang em right',
    'hang up basketball hoop',
    'happy knees',
    'hard back brace',
    'hard pastel sticks',
    'hate race',
    'heavy basketball training balls',
    'heavy bikes games',
    'heavy compression pants',
    'heavy training basketball',
    'heavy vest for running',
    'heavy weight basketball',
    'hed bands',
    'helping women recover',
    'hem pants tape',
    'hey black',
    'hey lady',
    'high and tight',
    'high compression',
    'high compression leggings',
    'high compression leggings women',
    'high compression pants women',
    'high compression women pants yoga fitness leggings',
    'high hip pain',
    'high intensity leggings',
    'high leggings',
    'high lift recovery',
    'high pants',
    'high pants for women',
    'high performance yoga pants',
    'high power capri',
    'high ride',
    'high rise bike shorts',
    'high rise bike shorts women',
    'high rise black leggings',
    'high rise capri',
    'high rise capri leggings',
    'high rise capri leggings for women',
    'high rise capris',
    'high rise capris athletic',
    'high rise capris for women',
    'high rise compression capri',
    'high rise compression leggings women',
    'high rise compression workout capris',
    'high rise dress pants',
    'high rise green leggings',
    'high rise legging',
    'high rise leggings',
    'high rise leggings for women',
    'high rise leggings workout',
    'high rise leggins',
    'high rise running leggings',
    'high rise running leggings for women',
    'high rise running shorts',
    'high rise running shorts women',
    'high rise shorts men',
    'high rise trousers',
    'high rise workout leggings for women',
    'high rise workout shorts',
    'high tech knee brace',
    'high tech workout',
    'high thighs',
    'high tights',
    'high top nike',
    'high top nike for men',
    'high top nike for women',
    'high top nikes',
    'high top puma',
    'high visibility clothing for women',
    'high visibility fabric tape',
    'high visibility running jacket',
    'high visibility shirts women',
    'high visibility waist running',
    'high visibility winter wear',
    'high waist',
    'high waist 90 degree',
    'high waist active wear',
    'high waist athletic leggings',
    'high waist bike shorts for women',
    'high waist black',
    'high waist black pants',
    'high waist boxing shorts',
    'high waist capri',
    'high waist capri leggings',
    'high waist capri leggings for women',
    'high waist capri pants',
    'high waist capri pants for women',
    'high waist capri tights',
    'high waist capri workout pants'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
